$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D14").Formula = "=(SUM(D3:D13)/6)"
$ws.Range("B13").Select()
